# Auto-generated edit script: adds the 'Knärot' section (with references) before
# the trailing sectPr, and updates the title-page date.

$d = $word.ActiveDocument

$paragraphData = @(
    @{ Style = 'Heading1'; Runs = @(
        @{ Italic = $false; Text = 'Knärot – ekologi samt krav på livsmiljön' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Samuel Johnsons doktorsavhandling ' },
        @{ Italic = $true; Text = '“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“' },
        @{ Italic = $false; Text = ' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ' },
        @{ Italic = $true; Text = '“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ' },
        @{ Italic = $false; Text = 'Vidare ' },
        @{ Italic = $true; Text = '“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ' },
        @{ Italic = $true; Text = '“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).' }
    ) },
    @{ Style = 'Heading2'; Runs = @(
        @{ Italic = $false; Text = 'Referenser - knärot' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'de Graaf M & Roberts M.R., 2009. ' },
        @{ Italic = $true; Text = 'Short-term response of the herbaceous layer within leave patches after harvest. ' },
        @{ Italic = $false; Text = 'Forest Ecology and Management 257, 1014-1025' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ' },
        @{ Italic = $true; Text = 'Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ' },
        @{ Italic = $false; Text = 'Ecological Applications, 22, 2049-2064 ' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ' },
        @{ Italic = $true; Text = 'Interactive effects of drought and edge exposure on old-growth forest understory species. ' },
        @{ Italic = $false; Text = 'Landscape Ecology, 37, sid 1839-1853' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ' },
        @{ Italic = $true; Text = 'Biological legacies buffer local species extinction after logging. ' },
        @{ Italic = $false; Text = 'Journal of Applied Ecology. 51, 53-62.' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'Skogsstyrelsen, 2022. ' },
        @{ Italic = $true; Text = 'Vägledning för hänsyn till knärot. ' },
        @{ Italic = $false; Text = 'https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/' }
    ) },
    @{ Style = 'Normal'; Runs = @(
        @{ Italic = $false; Text = 'SLU Artdatabanken, 2021. ' },
        @{ Italic = $true; Text = 'Artfaktablad. Naturvård – artfakta. ' },
        @{ Italic = $false; Text = 'SLU Artdatabanken, Uppsala ' }
    ) }
)

# --- Pass 1: create all the new (initially empty) paragraphs, in order, and set
#     paragraph styles only where they actually change vs. the previous paragraph
#     (mirrors what Word itself would record). ---
$insertAfterRange = $d.Paragraphs.Last.Range
$prevStyle = $d.Paragraphs.Last.Style.NameLocal
$paraIndices = @()

foreach ($pdata in $paragraphData) {
    $insertAfterRange.InsertParagraphAfter()
    $idx = $d.Paragraphs.Count
    $paraIndices += $idx
    if ($pdata.Style -ne $prevStyle) {
        $d.Paragraphs.Item($idx).Style = $pdata.Style
        $prevStyle = $pdata.Style
    }
    $insertAfterRange = $d.Paragraphs.Item($idx).Range
}

# --- Pass 2: fill in the run text (and italics) for each paragraph, tracking the
#     insertion cursor as an explicit position so runs never bleed across the
#     paragraph-mark boundary. ---
for ($i = 0; $i -lt $paragraphData.Count; $i++) {
    $pdata = $paragraphData[$i]
    $idx = $paraIndices[$i]
    $cursor = $d.Paragraphs.Item($idx).Range.Start
    foreach ($rdata in $pdata.Runs) {
        $r = $d.Range($cursor, $cursor)
        $r.InsertAfter($rdata.Text)
        if ($rdata.Italic) {
            $r.Font.Italic = $true
        }
        $cursor = $r.End
    }
}

# --- Update the title-page date (header) from 2023-09-13 to 2023-09-15 ---
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}
